# Actualización automática 2025-10-06 16:30:09
# Set PRESUPUESTO (column G) values to 0 for the affected rows on the
# "VENTA MENSUAL" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

$rows = @(3, 6, 16, 17, 19, 20, 21, 22, 24, 32)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 7).Value = 0
}
